$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.677.52'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '2.230.17'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''251.94'
$ws.Range('E5').Value = '  +8.14%  '
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('D7').Value = '''71.03'
$ws.Range('E7').Value = '  +1.33%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.571'
$ws.Range('E9').Value = '  +1.77%  '
$ws.Range('D10').Value = '''42.72'
$ws.Range('E10').Value = '  +19.96%  '
$ws.Range('E11').Value = '  -2.78%  '
$ws.Range('D12').Value = '''59.01'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('E14').Value = '  +2.85%  '
$ws.Range('D15').Value = '2.562.11'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '''14.94'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('E17').Value = '  -1.28%  '
$ws.Range('D18').Value = '2.232.48'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').Value = '41.622.60'
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').Value = '0.0₃0970'
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').Value = '''6.21'
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').Value = '''73.03'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').Value = '''2.27'
$ws.Range('E23').Value = '  +10.30%  '
$ws.Range('D24').Value = '''234.84'
$ws.Range('E24').Value = '  -1.10%  '
$ws.Range('D25').Value = '''3.88'
$ws.Range('E25').Value = '  +6.55%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = '''2.50'
$ws.Range('E27').Value = '  +6.04%  '
$ws.Range('D28').Value = '''10.41'
$ws.Range('E28').Value = '  +3.62%  '
$ws.Range('E29').Value = '  +1.45%  '
$ws.Range('D30').Value = '''171.56'
$ws.Range('E30').Value = '  +2.10%  '
$ws.Range('D31').Value = '''20.66'
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('E33').Value = '  -1.85%  '
$ws.Range('D34').Value = '''5.58'
$ws.Range('E34').Value = '  +2.75%  '
$ws.Range('D35').Value = '''0.0721'
$ws.Range('E35').Value = '  +0.32%  '
$ws.Range('D36').Value = '''26.83'
$ws.Range('E36').Value = '  +20.30%  '
$ws.Range('D37').Value = '''4.66'
$ws.Range('E37').Value = '  -2.56%  '
$ws.Range('E38').Value = '  +10.81%  '
$ws.Range('D39').Value = '''0.0289'
$ws.Range('E39').Value = '  +7.88%  '
$ws.Range('E40').Value = '  +2.11%  '
$ws.Range('D41').Value = '''69.30'
$ws.Range('E41').Value = '  +2.46%  '
$ws.Range('D42').Value = '''6.02'
$ws.Range('E42').Value = '  -0.73%  '
$ws.Range('D43').Value = '''12.00'
$ws.Range('E43').Value = '  +18.87%  '
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('D45').Value = '''0.208'
$ws.Range('E45').Value = '  +9.80%  '
$ws.Range('D46').Value = '''8.83'
$ws.Range('E46').Value = '  -2.93%  '
$ws.Range('E47').Value = '  +9.63%  '
$ws.Range('E48').Value = '  +1.28%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').Value = '''1.16'
$ws.Range('E50').Value = '  +7.16%  '
$ws.Range('E51').Value = '  +1.92%  '

$ws.Range('D5').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D50').Style = "Normal"
